{"js": "// Apply the diff: update the date line and the two-digit multiplication\n// problems throughout the document. Every source string in this document\n// is unique, so a direct search-and-replace per pair is safe and will not\n// cross-contaminate other cells.\nconst replacements = [\n  [\"2025-06-20 Friday\", \"2025-06-21 Saturday\"],\n  [\"47\u00d740=\", \"33\u00d753=\"],\n  [\"20\u00d727=\", \"22\u00d763=\"],\n  [\"48\u00d764=\", \"49\u00d745=\"],\n  [\"89\u00d785=\", \"89\u00d714=\"],\n  [\"73\u00d751=\", \"15\u00d776=\"],\n  [\"15\u00d785=\", \"33\u00d758=\"],\n  [\"89\u00d771=\", \"59\u00d797=\"],\n  [\"46\u00d752=\", \"74\u00d755=\"],\n  [\"86\u00d764=\", \"91\u00d719=\"],\n  [\"74\u00d758=\", \"46\u00d770=\"],\n  [\"96\u00d796=\", \"19\u00d749=\"],\n  [\"28\u00d793=\", \"67\u00d787=\"],\n  [\"25\u00d772=\", \"54\u00d723=\"],\n  [\"70\u00d783=\", \"58\u00d775=\"],\n  [\"71\u00d765=\", \"90\u00d785=\"],\n  [\"93\u00d742=\", \"99\u00d749=\"],\n  [\"73\u00d728=\", \"59\u00d729=\"],\n  [\"66\u00d747=\", \"43\u00d746=\"],\n  [\"12\u00d731=\", \"96\u00d799=\"],\n  [\"54\u00d798=\", \"11\u00d712=\"],\n  [\"60\u00d736=\", \"55\u00d781=\"],\n  [\"26\u00d769=\", \"99\u00d786=\"],\n  [\"94\u00d724=\", \"40\u00d716=\"],\n  [\"30\u00d793=\", \"29\u00d768=\"],\n  [\"91\u00d754=\", \"79\u00d755=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the diff: update the date line and the two-digit multiplication\n# problems throughout the document. Every source string in this document\n# is unique, so a direct Find/Replace per pair is safe and will not\n# cross-contaminate other cells.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-06-20 Friday\", \"2025-06-21 Saturday\"),\n    @(\"47\u00d740=\", \"33\u00d753=\"),\n    @(\"20\u00d727=\", \"22\u00d763=\"),\n    @(\"48\u00d764=\", \"49\u00d745=\"),\n    @(\"89\u00d785=\", \"89\u00d714=\"),\n    @(\"73\u00d751=\", \"15\u00d776=\"),\n    @(\"15\u00d785=\", \"33\u00d758=\"),\n    @(\"89\u00d771=\", \"59\u00d797=\"),\n    @(\"46\u00d752=\", \"74\u00d755=\"),\n    @(\"86\u00d764=\", \"91\u00d719=\"),\n    @(\"74\u00d758=\", \"46\u00d770=\"),\n    @(\"96\u00d796=\", \"19\u00d749=\"),\n    @(\"28\u00d793=\", \"67\u00d787=\"),\n    @(\"25\u00d772=\", \"54\u00d723=\"),\n    @(\"70\u00d783=\", \"58\u00d775=\"),\n    @(\"71\u00d765=\", \"90\u00d785=\"),\n    @(\"93\u00d742=\", \"99\u00d749=\"),\n    @(\"73\u00d728=\", \"59\u00d729=\"),\n    @(\"66\u00d747=\", \"43\u00d746=\"),\n    @(\"12\u00d731=\", \"96\u00d799=\"),\n    @(\"54\u00d798=\", \"11\u00d712=\"),\n    @(\"60\u00d736=\", \"55\u00d781=\"),\n    @(\"26\u00d769=\", \"99\u00d786=\"),\n    @(\"94\u00d724=\", \"40\u00d716=\"),\n    @(\"30\u00d793=\", \"29\u00d768=\"),\n    @(\"91\u00d754=\", \"79\u00d755=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute(\n        $oldText,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $newText,\n        2\n    )\n}\n"}
